# Applies the "integrate ramping into data prep and input data" edit:
#  - Definition sheet: reorder the node rows (A7:A18) into their new order
#  - Definition_parameters sheet: insert a new "unit_on_cost" row for
#    Methanol_Reactor, add the "online_variable_type" row, and shift the
#    power_line_Wholesale_Kasso row down
#  - Nodes sheet: reorder the node rows (2:13) into their new order, and
#    give Waste_Heat a node_slack_penalty of 100000
#  - Object__to_from_node sheet: add ramp_up_limit / ramp_down_limit rows
#    for Methanol_Reactor -> Raw_Methanol, and reorder / re-point the
#    remaining connection rows

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Definition" sheet - reorder A7:A18 (Category column B stays "node")
# ---------------------------------------------------------------------
$wsDef = $wb.Worksheets.Item("Definition")
$defOrder = @(
    "Power_Wholesale",
    "E-Methanol_storage_Kasso",
    "Vaporized_Carbon_Dioxide",
    "Power_Kasso",
    "Carbon_Dioxide",
    "Hydrogen_storage_Kasso",
    "Raw_Methanol",
    "District_Heating",
    "E-Methanol_Kasso",
    "Waste_Heat",
    "Water",
    "Hydrogen_Kasso"
)
for ($i = 0; $i -lt $defOrder.Length; $i++) {
    $row = 7 + $i
    $wsDef.Cells.Item($row, 1).Value = $defOrder[$i]
}

# ---------------------------------------------------------------------
# 2. "Definition_parameters" sheet - insert a row so the Methanol_Reactor
#    parameters become: unit_on_cost, min_down_time, online_variable_type
# ---------------------------------------------------------------------
$wsDefParam = $wb.Worksheets.Item("Definition_parameters")
# Insert a new row above the current row 3 (min_down_time), pushing
# min_down_time / online_variable_type / power_line_Wholesale_Kasso down.
$wsDefParam.Rows.Item(3).Insert()

# Row 3: Methanol_Reactor / unit / unit_on_cost / 1e-07  (new row)
$wsDefParam.Cells.Item(3, 1).Value = "Methanol_Reactor"
$wsDefParam.Cells.Item(3, 2).Value = "unit"
$wsDefParam.Cells.Item(3, 3).Value = "unit_on_cost"
$wsDefParam.Cells.Item(3, 4).Value = 0.0000001

# Row 4: Methanol_Reactor / unit / min_down_time / 48  (shifted down, values unchanged)
$wsDefParam.Cells.Item(4, 1).Value = "Methanol_Reactor"
$wsDefParam.Cells.Item(4, 2).Value = "unit"
$wsDefParam.Cells.Item(4, 3).Value = "min_down_time"
$wsDefParam.Cells.Item(4, 4).Value = 48

# Row 5: Methanol_Reactor / unit / online_variable_type / unit_online_variable_type_integer
$wsDefParam.Cells.Item(5, 1).Value = "Methanol_Reactor"
$wsDefParam.Cells.Item(5, 2).Value = "unit"
$wsDefParam.Cells.Item(5, 3).Value = "online_variable_type"
$wsDefParam.Cells.Item(5, 4).Value = "unit_online_variable_type_integer"

# Row 6: power_line_Wholesale_Kasso / connection / fom_cost / 100 (shifted down, values unchanged)
$wsDefParam.Cells.Item(6, 1).Value = "power_line_Wholesale_Kasso"
$wsDefParam.Cells.Item(6, 2).Value = "connection"
$wsDefParam.Cells.Item(6, 3).Value = "fom_cost"
$wsDefParam.Cells.Item(6, 4).Value = 100

# ---------------------------------------------------------------------
# 3. "Nodes" sheet - reorder rows 2:13, and give Waste_Heat a
#    node_slack_penalty (column G) of 100000
# ---------------------------------------------------------------------
$wsNodes = $wb.Worksheets.Item("Nodes")
# Name, balance_type, has_state, node_state_cap, frac_state_loss, node_slack_penalty
$nodesData = @(
    @("Power_Wholesale",          "balance_type_none", $null,  $null,   $null, $null),
    @("E-Methanol_storage_Kasso", "balance_type_node", "true", 100000,  0,     100000),
    @("Vaporized_Carbon_Dioxide", "balance_type_node", $null,  $null,   $null, 100000),
    @("Power_Kasso",              "balance_type_node", $null,  $null,   $null, 100000),
    @("Carbon_Dioxide",           "balance_type_none", $null,  $null,   $null, $null),
    @("Hydrogen_storage_Kasso",   "balance_type_node", "true", 100000,  0,     100000),
    @("Raw_Methanol",             "balance_type_node", $null,  $null,   $null, 100000),
    @("District_Heating",         "balance_type_none", $null,  $null,   $null, $null),
    @("E-Methanol_Kasso",         "balance_type_node", $null,  $null,   $null, 100000),
    @("Waste_Heat",               "balance_type_node", $null,  $null,   $null, 100000),
    @("Water",                    "balance_type_none", $null,  $null,   $null, $null),
    @("Hydrogen_Kasso",           "balance_type_node", $null,  $null,   $null, 100000)
)
for ($i = 0; $i -lt $nodesData.Length; $i++) {
    $row = 2 + $i
    $entry = $nodesData[$i]
    $wsNodes.Cells.Item($row, 1).Value = $entry[0]
    $wsNodes.Cells.Item($row, 2).Value = "node"
    $wsNodes.Cells.Item($row, 3).Value = $entry[1]
    if ($null -eq $entry[2]) { $wsNodes.Cells.Item($row, 4).Value = "" } else { $wsNodes.Cells.Item($row, 4).Value = "'" + $entry[2] }
    if ($null -eq $entry[3]) { $wsNodes.Cells.Item($row, 5).Value = "" } else { $wsNodes.Cells.Item($row, 5).Value = $entry[3] }
    if ($null -eq $entry[4]) { $wsNodes.Cells.Item($row, 6).Value = "" } else { $wsNodes.Cells.Item($row, 6).Value = $entry[4] }
    if ($null -eq $entry[5]) { $wsNodes.Cells.Item($row, 7).Value = "" } else { $wsNodes.Cells.Item($row, 7).Value = $entry[5] }
}

# ---------------------------------------------------------------------
# 4. "Object__to_from_node" sheet - rework rows 8-17 and append 2 new rows
#    for the Methanol_Reactor ramping parameters
# ---------------------------------------------------------------------
$wsObj = $wb.Worksheets.Item("Object__to_from_node")

# Insert 2 new rows after row 17 (before the final pipeline_District_Heating row)
$wsObj.Rows.Item(18).Insert()
$wsObj.Rows.Item(18).Insert()

# relationship_class_name, object_class, object_name, node, parameter_name, value
$objRows = @{
    8  = @("unit__to_node",           "unit",       "Methanol_Reactor",          "Raw_Methanol",             "ramp_up_limit",      0.5)
    9  = @("unit__to_node",           "unit",       "Methanol_Reactor",          "Raw_Methanol",             "ramp_down_limit",    0.5)
    10 = @("unit__to_node",           "unit",       "Methanol_Reactor",          "Waste_Heat",               "unit_capacity",      100)
    11 = @("connection__from_node",   "connection", "power_line_Wholesale_Kasso","Power_Wholesale",          "connection_capacity",1000)
    12 = @("connection__to_node",     "connection", "power_line_Wholesale_Kasso","Power_Kasso",              "connection_capacity",1000)
    13 = @("connection__from_node",   "connection", "power_line_Wholesale_Kasso","Power_Kasso",              "connection_capacity",1000)
    14 = @("connection__to_node",     "connection", "power_line_Wholesale_Kasso","Power_Wholesale",          "connection_capacity",1000)
    15 = @("connection__to_node",     "connection", "pipeline_storage_hydrogen", "Hydrogen_storage_Kasso",   "connection_capacity",1000)
    16 = @("connection__from_node",   "connection", "pipeline_storage_hydrogen", "Hydrogen_storage_Kasso",   "connection_capacity",1000)
    17 = @("connection__to_node",     "connection", "pipeline_storage_e-methanol","E-Methanol_storage_Kasso","connection_capacity",1000)
    18 = @("connection__from_node",   "connection", "pipeline_storage_e-methanol","E-Methanol_storage_Kasso","connection_capacity",1000)
    19 = @("connection__from_node",   "connection", "pipeline_District_Heating", "Waste_Heat",               "connection_capacity",1000)
    20 = @("connection__to_node",     "connection", "pipeline_District_Heating", "District_Heating",         "connection_capacity",1000)
}
foreach ($row in ($objRows.Keys | Sort-Object)) {
    $entry = $objRows[$row]
    $wsObj.Cells.Item($row, 1).Value = $entry[0]
    $wsObj.Cells.Item($row, 2).Value = $entry[1]
    $wsObj.Cells.Item($row, 3).Value = $entry[2]
    $wsObj.Cells.Item($row, 4).Value = $entry[3]
    $wsObj.Cells.Item($row, 5).Value = $entry[4]
    $wsObj.Cells.Item($row, 6).Value = $entry[5]
}
